$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 70
$ws.Range("F6").Value = 831
$ws.Range("F7").Value = 412
$ws.Range("F8").Value = 4687
$ws.Range("F9").Value = 4687
$ws.Range("F11").Value = 119
$ws.Range("F12").Value = 155
$ws.Range("F15").Value = 115
$ws.Range("F16").Value = 7455
$ws.Range("F18").Value = 126
$ws.Range("F21").Value = 519
$ws.Range("F22").Value = 1354
$ws.Range("F24").Value = 6283
$ws.Range("F25").Value = 1741
$ws.Range("F27").Value = 1989
$ws.Range("F28").Value = 6162
$ws.Range("F29").Value = 141
$ws.Range("F31").Value = 116
$ws.Range("F33").Value = 445
$ws.Range("F34").Value = 6394
$ws.Range("F36").Value = 206
$ws.Range("F37").Value = 96
$ws.Range("F41").Value = 2454
$ws.Range("F43").Value = 58
$ws.Range("F45").Value = 38
$ws.Range("F46").Value = 427
$ws.Range("F47").Value = 2136
$ws.Range("F48").Value = 42
$ws.Range("F49").Value = 1073

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 231
$ws.Range("F4").Value = 9
$ws.Range("F6").Value = 123
$ws.Range("F9").Value = 43
$ws.Range("F14").Value = 22

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1443

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1443
$ws.Range("F5").Value = 231
$ws.Range("F6").Value = 70
$ws.Range("F8").Value = 412
$ws.Range("F9").Value = 4687
$ws.Range("F10").Value = 4687
$ws.Range("F12").Value = 119
$ws.Range("F13").Value = 155
$ws.Range("F16").Value = 115
$ws.Range("F17").Value = 7455
$ws.Range("F19").Value = 126
$ws.Range("F20").Value = 519
$ws.Range("F21").Value = 1354
$ws.Range("F22").Value = 123
$ws.Range("F23").Value = 6283
$ws.Range("F24").Value = 1741
$ws.Range("F26").Value = 1989
$ws.Range("F28").Value = 43
$ws.Range("F29").Value = 6162
$ws.Range("F30").Value = 141
$ws.Range("F33").Value = 116
$ws.Range("F35").Value = 445
$ws.Range("F36").Value = 6394
$ws.Range("F38").Value = 206
$ws.Range("F39").Value = 96
$ws.Range("F42").Value = 2454
$ws.Range("F45").Value = 38
$ws.Range("F46").Value = 427
$ws.Range("F48").Value = 2136
$ws.Range("F49").Value = 42
$ws.Range("F50").Value = 22
